$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: "Channels and Guide seem redundant." -> mark Status as "Done-ish"
$ws.Range("B27").Value = "Done-ish"

# Row 31: "Guide should roll up and down with the mouse wheel..." -> mark Status as "Done"
# and update the Comments text to describe how it was done.
$ws.Range("B31").Value = "Done"
$ws.Range("C31").Value = "Done using a keymap. Ctrl-g for guide with support of a script."

# Update the view so it scrolls to show row 19 at top and the active selection is C23.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$ws.Range("C23").Select()
